$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nulos_por_campo")

$ws.Range("A2").Value = "tienda"
$ws.Range("B2").Value = 1200

$ws.Range("A3").Value = "satisfaccion"
$ws.Range("B3").Value = 227

$ws.Range("A4").Value = "comentario"
$ws.Range("B4").Value = 201

$ws.Range("A5").Value = "canal"
$ws.Range("B5").Value = 0

$ws.Range("A6").Value = "id_respuesta"
$ws.Range("B6").Value = 0

$ws.Range("A7").Value = "fecha"
$ws.Range("B7").Value = 0

$ws.Range("A8").Value = "producto"
$ws.Range("B8").Value = 0
